$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended by the Adafruit IO sync (row 53), mirroring the
# shape/content of the preceding rows in the sheet.
$row = 53

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C holds a numeric-looking value ("25"). Force it to be stored as
# literal text (matching the rest of the "Value" column in the sheet)
# instead of letting Excel auto-convert it to a number.
$valueCell = $ws.Cells.Item($row, 3)
$valueCell.NumberFormat = "@"
$valueCell.Value = "25"
$valueCell.Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
